$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Update the data table (FirstName / MiddleName / LastName)
# Order matters for the shared-strings table append order, so update
# row 4, then row 5, then row 3 (matching the original authoring order).
$ws.Range("A4").Value = "Donald"
$ws.Range("C4").Value = "Trump"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("C5").Value = "Salah"

$ws.Range("A3").Value = "Katie"
$ws.Range("C3").Value = "Ball"

# Move the active selection to B7 as in the saved file
$ws.Activate()
$ws.Range("B7").Select()
